# New crime data collected — weekly CompStat 061 precinct update
# Advance the report week (1/15/2024-1/21/2024 -> 1/22/2024-1/28/2024),
# bump the volume/number header, and refresh the Crime Complaints
# table (rows 15-27, columns C:N) with the newly collected figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text: "Volume 31   Number  3" -> "...Number  4" ---
$ws.Range("A8").Characters(21, 1).Text = "4"

# --- Report date range: 1/15/2024-1/21/2024 -> 1/22/2024-1/28/2024 ---
$ws.Range("C9").Characters(27, 9).Text = "1/22/2024"
$ws.Range("C9").Characters(47, 9).Text = "1/28/2024"

# --- Crime Complaints table (rows 15-27): updated weekly figures ---
$ws.Range("D15").Value = 1
$ws.Range("D15").NumberFormat = '#,##0'
$ws.Range("E15").Value = -100
$ws.Range("E15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("G15").Value = 2
$ws.Range("J15").Value = 2
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 66.666666666666
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = 18.75
$ws.Range("I16").Value = 19
$ws.Range("J16").Value = 16
$ws.Range("K16").Value = 18.75
$ws.Range("L16").Value = 72.727272727272
$ws.Range("M16").Value = 111.111111111111
$ws.Range("N16").Value = -74.666666666666
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 200
$ws.Range("F17").Value = 19
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 19
$ws.Range("J17").Value = 19
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 5.555555555555
$ws.Range("M17").Value = 137.5
$ws.Range("N17").Value = -29.629629629629
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 8
$ws.Range("E18").Value = -37.5
$ws.Range("F18").Value = 14
$ws.Range("G18").Value = 26
$ws.Range("H18").Value = -46.153846153846
$ws.Range("I18").Value = 14
$ws.Range("J18").Value = 26
$ws.Range("K18").Value = -46.153846153846
$ws.Range("L18").Value = -12.5
$ws.Range("M18").Value = -22.222222222222
$ws.Range("N18").Value = -92.929292929292
$ws.Range("C19").Value = 13
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = -18.75
$ws.Range("F19").Value = 46
$ws.Range("G19").Value = 63
$ws.Range("H19").Value = -26.984126984127
$ws.Range("I19").Value = 46
$ws.Range("J19").Value = 63
$ws.Range("K19").Value = -26.984126984127
$ws.Range("L19").Value = -11.538461538461
$ws.Range("M19").Value = 53.333333333333
$ws.Range("N19").Value = -60.344827586206
$ws.Range("C20").Value = 7
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 250
$ws.Range("F20").Value = 17
$ws.Range("H20").Value = 88.888888888888
$ws.Range("I20").Value = 17
$ws.Range("J20").Value = 9
$ws.Range("K20").Value = 88.888888888888
$ws.Range("L20").Value = 30.769230769230
$ws.Range("M20").Value = 30.769230769230
$ws.Range("N20").Value = -93.089430894308
$ws.Range("C21").Value = 33
$ws.Range("D21").Value = 31
$ws.Range("E21").Value = 6.451612903225
$ws.Range("F21").Value = 115
$ws.Range("G21").Value = 136
$ws.Range("H21").Value = -15.441176470588
$ws.Range("I21").Value = 115
$ws.Range("J21").Value = 136
$ws.Range("K21").Value = -15.441176470588
$ws.Range("L21").Value = 4.545454545454
$ws.Range("M21").Value = 45.569620253164
$ws.Range("N21").Value = -82.654600301659
$ws.Range("D22").Value = 1
$ws.Range("D22").NumberFormat = '#,##0'
$ws.Range("E22").Value = -100
$ws.Range("E22").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = -50
$ws.Range("J22").Value = 2
$ws.Range("K22").Value = -50
$ws.Range("C23").Value = 2
$ws.Range("F23").Value = 13
$ws.Range("G23").Value = 3
$ws.Range("H23").Value = 333.333333333333
$ws.Range("I23").Value = 13
$ws.Range("K23").Value = 333.333333333333
$ws.Range("L23").Value = 1200
$ws.Range("M23").Value = 1200
$ws.Range("C24").Value = 20
$ws.Range("D24").Value = 16
$ws.Range("E24").Value = 25
$ws.Range("F24").Value = 77
$ws.Range("G24").Value = 96
$ws.Range("H24").Value = -19.791666666666
$ws.Range("I24").Value = 77
$ws.Range("J24").Value = 96
$ws.Range("K24").Value = -19.791666666666
$ws.Range("L24").Value = -12.5
$ws.Range("M24").Value = -11.494252873563
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = 125
$ws.Range("G25").Value = 27
$ws.Range("H25").Value = 3.703703703703
$ws.Range("I25").Value = 28
$ws.Range("J25").Value = 27
$ws.Range("K25").Value = 3.703703703703
$ws.Range("L25").Value = 27.272727272727
$ws.Range("M25").Value = -9.677419354838
$ws.Range("D26").Value = 1
$ws.Range("D26").NumberFormat = '#,##0'
$ws.Range("E26").Value = -100
$ws.Range("E26").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("G26").Value = 2
$ws.Range("J26").Value = 2
$ws.Range("C27").Value = 3
$ws.Range("D27").Value = 3
$ws.Range("D27").NumberFormat = '#,##0'
$ws.Range("E27").Value = 0
$ws.Range("E27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F27").Value = 5
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 25
$ws.Range("I27").Value = 5
$ws.Range("J27").Value = 4
$ws.Range("K27").Value = 25
$ws.Range("L27").Value = 400
